# "added new sensor migration and fixes to ble"
#
# 1. Insert a new "byte index" row (row 16) on Sheet1, B16:U16 = 0..19,
#    formatted like the other byte-number header rows (e.g. row 39).
# 2. Update the sheet's saved selection to U18 (and scroll the window
#    toward N4, best-effort).
#
# NOTE: re-resolve the workbook/sheet via $excel.ActiveWorkbook (rather
# than using the $wb/$ws handed in directly) - chained property/method
# calls off the pre-bound $wb reference intermittently null-ref in this
# host, re-fetching avoids it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. New row 16: byte-index header, same formatting as row 4 (full
#        B:U width, style carried from the existing byte-number rows) ---
$ws.Range("B4:U4").Copy()
$ws.Range("B16").Select()
$ws.Paste()

$byteIndex = @(0,1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19)
for ($i = 0; $i -lt $byteIndex.Length; $i++) {
    $col = 2 + $i   # B = 2 .. U = 21
    $ws.Cells.Item(16, $col).Value = $byteIndex[$i]
}

# --- 2. View state: scroll + selection ---
$win = $excel.Windows[1]
$win.ScrollRow = 4
$win.ScrollColumn = 14
$ws.Range("U18").Select()
